$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (895 / Zika) entirely - shrinks dimension to A1:E35
$ws.Rows(36).Delete()

# Row 3 (113 - Desnutricion aguda)
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 0.27

# Row 5 (155 - Cancer de la mama y cuello uterino)
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 0.2

# Row 6 (210 - Dengue)
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 0.01

# Row 7 (215 - Defectos congenitos)
$ws.Range("D7").Value = 7
$ws.Range("E7").Value = 0.02

# Row 9 (300 - Agresiones por animales potencialmente transmisores de rabia)
$ws.Range("C9").Value = 40
$ws.Range("D9").Value = 47
$ws.Range("E9").Value = 0.03

# Row 10 (330 - Hepatitis a)
$ws.Range("C10").Value = 0
$ws.Range("E10").Value = 1

# Row 11 (340 - Hepatitis b, c y coinfeccion hepatitis b y delta)
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 0

# Row 12 (342 - Enfermedades huerfanas - raras)
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 0.22

# Row 13 (346 - Ira por virus nuevo)
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 0.18

# Row 16 (355 - Enfermedad transmitida por alimentos o agua (eta))
$ws.Range("C16").Value = 4
$ws.Range("E16").Value = 0.02

# Row 17 (356 - Intento de suicidio)
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 0.11

# Row 19 (365 - Intoxicaciones)
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 1

# Row 22 (455 - Leptospirosis)
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0.37

# Row 25 (549 - Morbilidad materna extrema)
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 0.09

# Row 26 (560 - Mortalidad perinatal y neonatal tardia)
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 0

# Row 29 (620 - Parotiditis)
$ws.Range("C29").Value = 2
$ws.Range("E29").Value = 0.14

# Row 31 (750 - Sifilis gestacional)
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0.14

# Row 33 (813 - Tuberculosis)
$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 5
$ws.Range("E33").Value = 0.16

# Row 34 (831 - Varicela individual)
$ws.Range("D34").Value = 3
$ws.Range("E34").Value = 0.01

# Row 35 (850 - Vih/sida/mortalidad por sida)
$ws.Range("D35").Value = 6
